# Applies the "add larger test dataset for admins" edit:
#  - Swaps the "Current status" of the first two sites in "Site data"
#    (Johnston's Pond: Maintain -> Restore, Lobster Bay: Restore -> Maintain)
#  - Updates the cost figures (columns E/F/G) on "Site data"
#  - Updates the feature goal/weight figures on "Feature data"
#  - Updates the expected feature amounts on the three "Consequence of ..." sheets

$wb = $excel.ActiveWorkbook

# NOTE: these worksheets have (password-less) sheet protection turned on,
# but all of the data-entry cells we touch are already unlocked. Using
# `.Value2` (rather than `.Value`) writes through without needing to
# Unprotect/Protect the sheet (which would also risk disturbing the
# <sheetProtection> element that the diff leaves untouched).

# ---- Sheet 1: "Site data" ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("D4").Value2 = "Restore"
$ws1.Range("E4").Value2 = 151
$ws1.Range("F4").Value2 = 270
$ws1.Range("G4").Value2 = 487

$ws1.Range("D5").Value2 = "Maintain"
$ws1.Range("E5").Value2 = 171
$ws1.Range("F5").Value2 = 276
$ws1.Range("G5").Value2 = 526

$ws1.Range("D6").Value2 = "Maintain"
$ws1.Range("E6").Value2 = 170
$ws1.Range("F6").Value2 = 268
$ws1.Range("G6").Value2 = 442

$ws1.Range("D7").Value2 = "Restore"
$ws1.Range("E7").Value2 = 125
$ws1.Range("F7").Value2 = 290
$ws1.Range("G7").Value2 = 391

# ---- Sheet 3: "Feature data" ----
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B4").Value2 = 66
$ws3.Range("C4").Value2 = 32

$ws3.Range("B5").Value2 = 7
$ws3.Range("C5").Value2 = 22

$ws3.Range("B6").Value2 = 64
$ws3.Range("C6").Value2 = 78

# ---- Sheet 4: "Consequence of "Maintain"" ----
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B4").Value2 = 141
$ws4.Range("C4").Value2 = 178
$ws4.Range("D4").Value2 = 113

$ws4.Range("B5").Value2 = 107
$ws4.Range("C5").Value2 = 144
$ws4.Range("D5").Value2 = 173

$ws4.Range("B6").Value2 = 197
$ws4.Range("C6").Value2 = 177
$ws4.Range("D6").Value2 = 151

$ws4.Range("B7").Value2 = 165
$ws4.Range("C7").Value2 = 175
$ws4.Range("D7").Value2 = 131

# ---- Sheet 5: "Consequence of "Signage"" ----
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B4").Value2 = 204
$ws5.Range("C4").Value2 = 275
$ws5.Range("D4").Value2 = 280

$ws5.Range("B5").Value2 = 377
$ws5.Range("C5").Value2 = 392
$ws5.Range("D5").Value2 = 328

$ws5.Range("B6").Value2 = 341
$ws5.Range("C6").Value2 = 334
$ws5.Range("D6").Value2 = 245

$ws5.Range("B7").Value2 = 388
$ws5.Range("C7").Value2 = 279
$ws5.Range("D7").Value2 = 297

# ---- Sheet 6: "Consequence of "Restore"" ----
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B4").Value2 = 573
$ws6.Range("C4").Value2 = 485
$ws6.Range("D4").Value2 = 366

$ws6.Range("B5").Value2 = 338
$ws6.Range("C5").Value2 = 524
$ws6.Range("D5").Value2 = 426

$ws6.Range("B6").Value2 = 423
$ws6.Range("C6").Value2 = 358
$ws6.Range("D6").Value2 = 540

$ws6.Range("B7").Value2 = 337
$ws6.Range("C7").Value2 = 468
$ws6.Range("D7").Value2 = 428

# ---- Sheet 7: "metadata" (hidden) ----
# Action labels for the first/third action stay the same text values
# ("Maintain" / "Restore") even though the underlying shared-string
# table gets reshuffled by the site-status swap above.
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("B2").Value2 = "Maintain"
$ws7.Range("B4").Value2 = "Restore"
